$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.628069043159485
$ws.Range("B1").Value = 4.069114208221436
$ws.Range("C1").Value = 2.766246318817139
$ws.Range("D1").Value = 0.9726204872131348
$ws.Range("E1").Value = 0.8799508213996887
